# Apply the "B suite" test-case sheet update:
#  - fix column D style on rows 2-45 (border-only, no fill)
#  - row 45 Results goes back from PASS to SKIP (it is no longer the last row)
#  - append two new test cases as rows 46 and 47
#  - refresh the sheet selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---------------------------------------------------------------------------
# 1. Normalize column D (Runmode) styling for the existing data rows.
#    (Border-only cell format, matching the rest of the sheet.)
# ---------------------------------------------------------------------------
$ws.Range("D2:D45").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2. Row 45 no longer is the final row, so its Results value reverts to SKIP.
# ---------------------------------------------------------------------------
$ws.Range("E45").Value = "SKIP"

# ---------------------------------------------------------------------------
# 3. Add the two new rows of test case data.
# ---------------------------------------------------------------------------
$ws.Range("A46").Value = "TestCase_B45"
$ws.Range("B46").Value = "OPQA-270"
$ws.Range("C46").Value = "Verify that following  content type options are present in the search drop down:`na)All`nb)Articles`nc)Patents`nd)People`ne)Posts"
$ws.Range("D46").Value = "Y"
$ws.Range("E46").Value = "SKIP"

$ws.Range("A47").Value = "TestCase_B46"
$ws.Range("B47").Value = "OPQA-274"
$ws.Range("C47").Value = "Verify that ALL content type is selected in the search drop down by default"
$ws.Range("D47").Value = "Y"
$ws.Range("E47").Value = "PASS"

# ---------------------------------------------------------------------------
# 4. Match formatting of the new rows to the rest of the table by copying
#    the format from equivalent existing cells.
# ---------------------------------------------------------------------------
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A46:A47").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B46:B47").PasteSpecial(-4122) | Out-Null

$ws.Range("C37").Copy() | Out-Null
$ws.Range("C46:C47").PasteSpecial(-4122) | Out-Null

$ws.Range("D45").Copy() | Out-Null
$ws.Range("D46:D47").PasteSpecial(-4122) | Out-Null

$ws.Range("E45").Copy() | Out-Null
$ws.Range("E46:E47").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Row 46 wraps 6 lines of text - match the row height used elsewhere on the
# sheet for similar multi-line entries.
$ws.Rows.Item(46).RowHeight = 90

# ---------------------------------------------------------------------------
# 5. Update the view: selection now spans the new rows and the window has
#    scrolled down so row 42 is the first visible row.
# ---------------------------------------------------------------------------
$ws.Range("A42").Activate() | Out-Null
$ws.Range("D2:D47").Select() | Out-Null
$ws.Range("D47").Activate() | Out-Null
